$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 9797.857
$ws.Range("I7").Value = 4986
$ws.Range("J7").Value = 10599.833
$ws.Range("K7").Value = 4986
$ws.Range("L7").Value = 10599.833
$ws.Range("M7").Value = -4874
$ws.Range("N7").Value = -10823.833

$ws.Range("H14").Value = 9797.857
$ws.Range("I14").Value = 4986
$ws.Range("J14").Value = 10599.833
$ws.Range("K14").Value = 4986
$ws.Range("L14").Value = 10599.833
$ws.Range("M14").Value = -4795
$ws.Range("N14").Value = -10981.833

$ws.Range("H113").Value = 6949.5
$ws.Range("I113").Value = 6949.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 6949.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -3695.5
$ws.Range("N113").ClearContents()

$ws.Range("H116").Value = 7315.625
$ws.Range("I116").Value = 5058.3335
$ws.Range("K116").Value = 5058.3335
$ws.Range("M116").Value = -1616.3335

$ws.Range("H132").Value = 1652.4
$ws.Range("I132").Value = 1430.1786
$ws.Range("J132").Value = 2541.2856
$ws.Range("K132").Value = 4290.5358
$ws.Range("L132").Value = 7623.8568
$ws.Range("M132").Value = -1760.5358
$ws.Range("N132").Value = -12683.8568

$ws.Range("H141").Value = 1616.591
$ws.Range("I141").Value = 1616.591
$ws.Range("K141").Value = 4849.772999999999
$ws.Range("M141").Value = 330.2270000000008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1078.25
$ws.Range("I5").Value = 1473.5
$ws.Range("K5").Value = 1473.5
$ws.Range("M5").Value = -1361.5

$ws.Range("H21").Value = 13099.8
$ws.Range("J21").Value = 14124.75
$ws.Range("L21").Value = 14124.75
$ws.Range("N21").Value = -14872.75

$ws.Range("H61").Value = 5178.857
$ws.Range("I61").Value = 3816.8865
$ws.Range("K61").Value = 3816.8865
$ws.Range("M61").Value = -3604.8865

$ws.Range("H74").Value = 3765.1914
$ws.Range("I74").Value = 3130.4194
$ws.Range("K74").Value = 3130.4194
$ws.Range("M74").Value = -2256.4194

$ws.Range("H77").Value = 3765.1914
$ws.Range("I77").Value = 3130.4194
$ws.Range("K77").Value = 15652.097
$ws.Range("M77").Value = -11284.097

$ws.Range("H132").Value = 5611.45
$ws.Range("I132").Value = 7264.125
$ws.Range("J132").Value = 4509.6665
$ws.Range("K132").Value = 21792.375
$ws.Range("L132").Value = 13528.9995
$ws.Range("M132").Value = -19262.375
$ws.Range("N132").Value = -18588.9995

$ws.Range("H136").Value = 5178.857
$ws.Range("I136").Value = 3816.8865
$ws.Range("K136").Value = 11450.6595
$ws.Range("M136").Value = -8900.6595

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1078.25
$ws.Range("I4").Value = 1473.5
$ws.Range("K4").Value = 1473.5
$ws.Range("M4").Value = -1358.5

$ws.Range("H35").Value = 67995
$ws.Range("J35").Value = 67995
$ws.Range("L35").Value = 67995
$ws.Range("N35").Value = -68615

$ws.Range("H80").Value = 348.66666
$ws.Range("I80").Value = 232.83333
$ws.Range("J80").Value = 406.58334
$ws.Range("K80").Value = 232.83333
$ws.Range("L80").Value = 406.58334
$ws.Range("M80").Value = 765.1666700000001
$ws.Range("N80").Value = -2402.58334

$ws.Range("H83").Value = 348.66666
$ws.Range("I83").Value = 232.83333
$ws.Range("J83").Value = 406.58334
$ws.Range("K83").Value = 1164.16665
$ws.Range("L83").Value = 2032.9167
$ws.Range("M83").Value = 3827.83335
$ws.Range("N83").Value = -12016.9167

$ws.Range("H105").Value = 4154.9414
$ws.Range("I105").Value = 4395.357
$ws.Range("K105").Value = 4395.357
$ws.Range("M105").Value = -2648.357

$ws.Range("H132").Value = 69760
$ws.Range("J132").Value = 69760
$ws.Range("L132").Value = 69760
$ws.Range("N132").Value = -79880

$ws.Range("H134").Value = 6582.478
$ws.Range("I134").Value = 4227.475
$ws.Range("K134").Value = 12682.425
$ws.Range("M134").Value = -10147.425

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 266.07144
$ws.Range("I7").Value = 193.2
$ws.Range("J7").Value = 448.25
$ws.Range("K7").Value = 193.2
$ws.Range("L7").Value = 448.25
$ws.Range("M7").Value = -80.19999999999999
$ws.Range("N7").Value = -674.25

$ws.Range("H31").Value = 3025
$ws.Range("I31").Value = 1907.3334
$ws.Range("J31").Value = 3293.24
$ws.Range("K31").Value = 1907.3334
$ws.Range("L31").Value = 3293.24
$ws.Range("M31").Value = -1612.3334
$ws.Range("N31").Value = -3883.24

$ws.Range("H34").Value = 3025
$ws.Range("I34").Value = 1907.3334
$ws.Range("J34").Value = 3293.24
$ws.Range("K34").Value = 1907.3334
$ws.Range("L34").Value = 3293.24
$ws.Range("M34").Value = -1705.3334
$ws.Range("N34").Value = -3697.24

$ws.Range("H58").Value = 3285
$ws.Range("I58").Value = 1257.3226
$ws.Range("J58").Value = 10269.223
$ws.Range("K58").Value = 1257.3226
$ws.Range("L58").Value = 10269.223
$ws.Range("M58").Value = -1054.3226
$ws.Range("N58").Value = -10675.223

$ws.Range("H136").Value = 3285
$ws.Range("I136").Value = 1257.3226
$ws.Range("J136").Value = 10269.223
$ws.Range("K136").Value = 3771.9678
$ws.Range("L136").Value = 30807.669
$ws.Range("M136").Value = -1221.9678
$ws.Range("N136").Value = -35907.669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 41.3
$ws.Range("I2").Value = 39.73077
$ws.Range("J2").Value = 51.5
$ws.Range("K2").Value = 238.38462
$ws.Range("L2").Value = 309
$ws.Range("M2").Value = -125.38462
$ws.Range("N2").Value = -535

$ws.Range("H5").Value = 1530.3529
$ws.Range("I5").Value = 964
$ws.Range("J5").Value = 1766.3334
$ws.Range("K5").Value = 2892
$ws.Range("L5").Value = 5299.0002
$ws.Range("M5").Value = -2780
$ws.Range("N5").Value = -5523.0002

$ws.Range("H17").Value = 491.07693
$ws.Range("I17").Value = 600
$ws.Range("J17").Value = 423
$ws.Range("K17").Value = 1800
$ws.Range("L17").Value = 1269
$ws.Range("M17").Value = -1631
$ws.Range("N17").Value = -1607

$ws.Range("H34").Value = 4447.7
$ws.Range("J34").Value = 6285.2856
$ws.Range("L34").Value = 18855.8568
$ws.Range("N34").Value = -19023.8568

$ws.Range("H38").Value = 2648.7693
$ws.Range("I38").Value = 475.55554
$ws.Range("K38").Value = 1426.66662
$ws.Range("M38").Value = -1079.66662

$ws.Range("H39").Value = 5978.5713
$ws.Range("J39").Value = 7800
$ws.Range("L39").Value = 23400
$ws.Range("N39").Value = -23988

$ws.Range("H131").Value = 1324.475
$ws.Range("J131").Value = 1371
$ws.Range("L131").Value = 4113
$ws.Range("N131").Value = -14193

$ws.Range("H135").Value = 1530.3529
$ws.Range("I135").Value = 964
$ws.Range("J135").Value = 1766.3334
$ws.Range("K135").Value = 8676
$ws.Range("L135").Value = 15897.0006
$ws.Range("M135").Value = -6141
$ws.Range("N135").Value = -20967.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 21365.637
$ws.Range("J26").Value = 21365.637
$ws.Range("L26").Value = 21365.637
$ws.Range("N26").Value = -21925.637

$ws.Range("H50").Value = 21365.637
$ws.Range("J50").Value = 21365.637
$ws.Range("L50").Value = 21365.637
$ws.Range("N50").Value = -22361.637

$ws.Range("H102").Value = 5253
$ws.Range("I102").Value = 5670.6665
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 5670.6665
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -4048.6665
$ws.Range("N102").Value = -7244

$ws.Range("H122").Value = 4392.75
$ws.Range("I122").Value = 4690.3335
$ws.Range("K122").Value = 14071.0005
$ws.Range("M122").Value = -11621.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 910.43475
$ws.Range("I16").Value = 866.75
$ws.Range("K16").Value = 866.75
$ws.Range("M16").Value = -696.75

$ws.Range("H46").Value = 1177.7778
$ws.Range("I46").Value = 796.4
$ws.Range("K46").Value = 796.4
$ws.Range("M46").Value = -608.4

$ws.Range("H68").Value = 3734.9285
$ws.Range("I68").Value = 2254.889
$ws.Range("J68").Value = 6399
$ws.Range("K68").Value = 2254.889
$ws.Range("L68").Value = 6399
$ws.Range("M68").Value = -1505.889
$ws.Range("N68").Value = -7897

$ws.Range("H71").Value = 3734.9285
$ws.Range("I71").Value = 2254.889
$ws.Range("J71").Value = 6399
$ws.Range("K71").Value = 11274.445
$ws.Range("L71").Value = 31995
$ws.Range("M71").Value = -7530.445
$ws.Range("N71").Value = -39483

$ws.Range("H136").Value = 5102.725
$ws.Range("I136").Value = 4829.6
$ws.Range("K136").Value = 14488.8
$ws.Range("M136").Value = -11938.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H136").Value = 956.80646
$ws.Range("I136").Value = 893.2857
$ws.Range("K136").Value = 2679.8571
$ws.Range("M136").Value = -129.8571000000002
